$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (it only held the "5816812 - João Paulo Alves Silva" name in B/C,
# with no label in column A). Deleting it shifts rows 14-24 up by one, which
# realigns all the column-A labels from "Programa resumido:" down to "Requisitos:".
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) now needs to show the docente's name instead of the old
# "Complementar a formação..." objectives text.
$ws.Range("B10").Value2 = "5816812 - João Paulo Alves Silva"
$ws.Range("C10").Value2 = "5816812 - João Paulo Alves Silva"

# Row 13 (now "Programa resumido:") gets "Semestral".
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

# Row 15 (now "Programa:") gets the date "01/01/2012". Copy it from B8/C8 (which
# already stores this exact text) so it is kept as text instead of being
# auto-converted into a date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (now "Método:") gets the docente's name.
$ws.Range("B18").Value2 = "5816812 - João Paulo Alves Silva"
$ws.Range("C18").Value2 = "5816812 - João Paulo Alves Silva"

# Row 19 (now "Critério:") gets the old "Método" description text.
$ws.Range("B19").Value2 = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value2 = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# Row 20 (now "Norma de recuperação:") gets the old "Critério" text.
$ws.Range("B20").Value2 = "Provas e trabalhos."
$ws.Range("C20").Value2 = "Provas e trabalhos."

# Row 21 (now "Bibliografia:") gets the old "Norma de recuperação" text.
$ws.Range("B21").Value2 = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value2 = "Prova única com nota maior ou igual a 5,0 (cinco)."
